# Area_区域表.xlsx — add "boolean/isShape" column (D) to the config table,
# and rename the "range/范围" column (C) to "points/点集".
#
# Table layout (sheet "工作表1"):
#   Row1: type row       (int | string | int[][] | boolean)
#   Row2: field-name row (id | name | points | isShape)
#   Row3: Chinese labels  (区域 ID | 名称 | 点集 | 是否构成形状)
#   Row4: "Language" sub-header (column B)
#   Row5/6: example data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: type declarations -------------------------------------------
$ws.Range("D1").Value = "boolean"

# --- Row 2: English field names -----------------------------------------
$ws.Range("C2").Value = "points"
$ws.Range("D2").Value = "isShape"

# --- Row 3: Chinese field labels -----------------------------------------
$ws.Range("C3").Value = "点集"
$ws.Range("D3").Value = "是否构成形状"
